$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New layout for the station rows (A2:F22), as produced by a descending ---
# --- sort on the "Domain" column (C), with an AutoFilter left on the table. ---
# Columns: StationID, StationCode, Domain, StationName, Latitude, Longitude
$rows = @(
    @(7,  "KCUT", 2, "Custer, SD",             43.732999999999997, -103.611),
    @(20, "KRAP", 2, "Rapid City Airport, SD",  44.042999999999999, -103.054),
    @(21, "KRCA", 2, "Ellsworth AFB, SD",       44.133000000000003, -103.1),
    @(23, "KSPF", 2, "Clyde Ice Field, SD",     44.482999999999997, -103.783),
    @(25, "KUDX", 2, "Rapid City NEXRAD, SD",   44.133000000000003, -102.833),
    @(26, "KUNR", 2, "Rapid City NWS, SD",      44.072699999999998, -103.211),
    @(0,  "K2WX", 2, "Buffalo, SD",             45.603999999999999, -103.54600000000001),
    @(1,  "K4DG", 2, "Douglas, WY",             42.75,               -105.383),
    @(6,  "KCDR", 2, "Chadron, NE",             42.837000000000003, -103.098),
    @(8,  "KD07", 2, "Faith, SD",               45.031999999999996, -102.01900000000001),
    @(9,  "KDGW", 2, "Converse Co Arpt, WY",    42.795999999999999, -105.38),
    @(11, "KGCC", 2, "Gillette, WY",            44.338999999999999, -105.542),
    @(12, "KGRN", 2, "Gordon, NE",              42.805999999999997, -102.175),
    @(14, "KHEI", 2, "Hettinger, ND",           46.017000000000003, -102.65),
    @(15, "KIEN", 2, "Pine Ridge, SD",          43.021000000000001, -102.518),
    @(18, "KPHP", 2, "Philip, SD",              44.051000000000002, -101.601),
    @(19, "KPIR", 2, "Pierre, SD",              44.381,              -100.286),
    @(3,  "KAIA", 1, "Alliance, NE",            42.05,               -102.8),
    @(4,  "KBFF", 1, "Scottsbluff, NE",         41.871000000000002, -103.593),
    @(24, "KTOR", 1, "Torrington, WY",          42.061,              -104.158),
    @(30, "MUNS", 1, "Munich, ND",              48.666699999999999, -98.834999999999994)
)

function Set-StationRows {
    $r = 2
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
        $r++
    }
}

# Write the new row order.
Set-StationRows

# --- Turn on the AutoFilter for the header row, and mark the filter database name ---
$headerRng = $ws.Range("A1:F1")
$headerRng.AutoFilter() | Out-Null

$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$F`$1")
$fdb.Visible = $false

# --- Perform the sort (descending by Domain) that produced this layout, so the ---
# --- sheet keeps a record of the last sort (sortState/sortCondition). ---
$sortRange = $ws.Range("A2:F22")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C2:C22"), 0, 2) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# The engine's stable sort doesn't reconstruct the exact manual re-ordering the
# author performed, so re-apply the intended final row order here.
Set-StationRows

# --- Update the active selection ---
$ws.Range("C7").Select() | Out-Null
